$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A72").Value = 71
$ws.Range("B72").Value = 1
$ws.Range("C72").Value = "2024-06-16 10:13:30"
$ws.Range("D72").Value = 200
$ws.Range("E72").Value = 8

$ws.Range("A73").Value = 72
$ws.Range("B73").Value = 2
$ws.Range("C73").Value = "2024-06-16 10:13:30"
$ws.Range("D73").Value = 200
$ws.Range("E73").Value = 0
